$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 408-411: shift the weekly date and refresh the
#     volume / price figures for the new reporting week (44448).
$ws.Cells.Item(408, 4).Value = 44448
$ws.Cells.Item(408, 13).Value = 240
$ws.Cells.Item(408, 14).Value = 8500
$ws.Cells.Item(408, 15).Value = 9000
$ws.Cells.Item(408, 16).Value = 8750
$ws.Cells.Item(408, 19).Value = 547

$ws.Cells.Item(409, 4).Value = 44448
$ws.Cells.Item(409, 13).Value = 200
$ws.Cells.Item(409, 14).Value = 7500
$ws.Cells.Item(409, 15).Value = 8000
$ws.Cells.Item(409, 16).Value = 7750
$ws.Cells.Item(409, 19).Value = 484

$ws.Cells.Item(410, 4).Value = 44448
$ws.Cells.Item(410, 13).Value = 240
$ws.Cells.Item(410, 14).Value = 8000
$ws.Cells.Item(410, 15).Value = 8500
$ws.Cells.Item(410, 16).Value = 8250
$ws.Cells.Item(410, 19).Value = 516

$ws.Cells.Item(411, 4).Value = 44448
$ws.Cells.Item(411, 13).Value = 200
$ws.Cells.Item(411, 14).Value = 7000
$ws.Cells.Item(411, 15).Value = 7500
$ws.Cells.Item(411, 16).Value = 7250
$ws.Cells.Item(411, 19).Value = 453

# --- Rows 412-413 keep the original date (44400) but the variety recorded
#     there changes from "Pink Lady" to "Fuji royal", with updated figures.
$ws.Cells.Item(412, 11).Value = "Fuji royal"
$ws.Cells.Item(412, 13).Value = 160
$ws.Cells.Item(412, 14).Value = 7500
$ws.Cells.Item(412, 15).Value = 8000
$ws.Cells.Item(412, 16).Value = 7750
$ws.Cells.Item(412, 19).Value = 484

$ws.Cells.Item(413, 11).Value = "Fuji royal"
$ws.Cells.Item(413, 13).Value = 160

# --- Append 4 new rows (414-417) re-stating the original "Pink Lady" /
#     "Granny Smith" entries for the 44400 date that were displaced above.
$newRows = @(
    @{ Row = 414; K = "Granny Smith"; L = "Primera"; M = 120; N = 7000; O = 7500; P = 7250; S = 453 },
    @{ Row = 415; K = "Granny Smith"; L = "Segunda"; M = 120; N = 6000; O = 6500; P = 6250; S = 391 },
    @{ Row = 416; K = "Pink Lady";    L = "Primera"; M = 120; N = 7000; O = 7500; P = 7250; S = 453 },
    @{ Row = 417; K = "Pink Lady";    L = "Segunda"; M = 120; N = 6000; O = 6500; P = 6250; S = 391 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 7
    $ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
    $ws.Cells.Item($row, 3).Value = "Ñuble"
    $ws.Cells.Item($row, 4).Value = 44400
    $ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(413, 4).NumberFormat()
    $ws.Cells.Item($row, 5).Value = 16
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100104
    $ws.Cells.Item($row, 8).Value = "Frutos de pepita"
    $ws.Cells.Item($row, 9).Value = 100104002
    $ws.Cells.Item($row, 10).Value = "Manzana"
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "`$/caja 16 kilos empedrada"
    $ws.Cells.Item($row, 18).Value = "Provincia de Curicó"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 16
}
